# Add new vocabulary rows, formulas, formatting and an AutoFilter to the
# word-list sheet, matching the target edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: push the existing word ("さいけいせい"/"再形成") down to
#        row 4 by inserting two new rows above it. ---------------------------
$ws.Rows("2:3").Insert()

# --- 2. Header row extra (blank, but styled) cells C1:D1 --------------------
$ws.Range("C1:D1").Value = ""

# --- 3. New vocabulary words (rows 2,3 then 5-9; row 4 already holds the
#        original さいけいせい/再形成 pair). Column order follows the Japanese
#        reading (hiragana) alphabetical sort that the author applied. -------
$words = @(
    @{ row = 2; a = "いただきます"; b = "いただきます" },
    @{ row = 3; a = "きょうかしょ"; b = "教科書" },
    @{ row = 5; a = "さいしゅうび"; b = "最終日" },
    @{ row = 6; a = "さくらんぼう"; b = "桜ん坊" },
    @{ row = 7; a = "しんぱいごと"; b = "心配事" },
    @{ row = 8; a = "にんじゃがく"; b = "忍者学" },
    @{ row = 9; a = "ほうせきばこ"; b = "宝石箱" }
)

foreach ($w in $words) {
    $r = $w.row
    $ws.Cells.Item($r, 1).Value = $w.a
    $ws.Cells.Item($r, 2).Value = $w.b
}

# --- 4. Formulas for every data row 2-9: C = A&B concatenation, D = LEN(A) --
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 3).Formula = "=A$r&B$r"
    $ws.Cells.Item($r, 4).Formula = "=LEN(A$r)"
}

# --- 5. Fonts / styles ------------------------------------------------------
# Original data font (row 4, A4:B4) grows from 11pt to 12pt, keeping its
# (colourless) 游ゴシック font.
$ws.Range("A4:B4").Font.Size = 12

# New word rows (everything except row 4) use 12pt 游ゴシック with the theme
# text colour, vertically centred and wrapped.
$newWordRows = @(2, 3, 5, 6, 7, 8, 9)
foreach ($r in $newWordRows) {
    $rng = $ws.Range("A$r`:B$r")
    $rng.Font.Size = 12
    $rng.Font.Name = "游ゴシック"
    $rng.Font.ThemeColor = 1
    $rng.VerticalAlignment = -4108
    $rng.WrapText = $true
}

# Formula columns C & D (rows 1-9) share the plain 12pt coloured 游ゴシック
# font used above, without the wrap/centre alignment.
$ws.Range("C1:D9").Font.Size = 12
$ws.Range("C1:D9").Font.Name = "游ゴシック"
$ws.Range("C1:D9").Font.ThemeColor = 1

# --- 6. Leftover blank template cells below the table (rows 11-16) ---------
$tail = $ws.Range("A11:D16")
$tail.Value = "_tmp_"
$tail.ClearContents()

# --- 7. Column widths / row heights -----------------------------------------
$ws.Range("A1:B1").ColumnWidth = 15
$ws.Range("C1").ColumnWidth = 29.453125
$ws.Range("1:9").RowHeight = 20

# --- 8. AutoFilter over the table + the hidden _FilterDatabase name --------
$ws.Range("A1:D10").AutoFilter()
$ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$D`$10")
$ws.Names.Item("_xlnm._FilterDatabase").Visible = $false

# --- 9. Selection / view ------------------------------------------------
$ws.Range("I13").Select()

$wb.Save()
